$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.578119397163391
$ws.Range("B1").Value = 4.337272644042969
$ws.Range("C1").Value = 3.471343994140625
$ws.Range("D1").Value = 1.555417656898499
$ws.Range("E1").Value = 0.5726577043533325
